$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.872.50'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.808.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '702.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +11.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.05'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.807.01'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.68'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +13.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.18'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.449.38'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.808.59'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.851.55'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.70'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.30'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +18.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '479.44'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.38%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.88'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.35'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.43'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.959.55'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +15.59%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.50'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.54'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.57%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.55%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.48%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.758.43'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.52'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +6.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000335'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +24.19%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +12.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.969'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.46'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '160.33'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.79'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.41'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.18%  '
